$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.854.08"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.37%  "

$ws.Range("D3").Value = "'2.316.09"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.11%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'97.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.81%  "

$ws.Range("D6").Value = "'272.39"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.47%  "

$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("D9").Value = "'0.626"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.19%  "

$ws.Range("D10").Value = "'45.31"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.83%  "

$ws.Range("E11").Value = "  -1.58%  "

$ws.Range("D12").Value = "'8.01"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.64%  "

$ws.Range("E13").Value = "  +0.15%  "

$ws.Range("D14").Value = "'2.654.89"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.83%  "

$ws.Range("D15").Value = "'15.49"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.70%  "

$ws.Range("D16").Value = "'0.873"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +7.05%  "

$ws.Range("D17").Value = "'2.318.62"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.12%  "

$ws.Range("D18").Value = "'43.780.21"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.42%  "

$ws.Range("D19").Value = "'0.0000109"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.11%  "

$ws.Range("D20").Value = "'6.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.43%  "

$ws.Range("D21").Value = "'73.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.49%  "

$ws.Range("D22").Value = "'239.76"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.87%  "

$ws.Range("E23").Value = "  -1.64%  "

$ws.Range("D24").Value = "'9.45"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.71%  "

$ws.Range("E25").Value = "  +0.00%  "

$ws.Range("D26").Value = "'2.54"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.73%  "

$ws.Range("E27").Value = "  -0.65%  "

$ws.Range("D28").Value = "'3.50"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.74%  "

$ws.Range("E29").Value = "  +2.53%  "

$ws.Range("D30").Value = "'38.12"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.43%  "

$ws.Range("D31").Value = "'22.38"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.23%  "

$ws.Range("D32").Value = "'175.21"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.36%  "

$ws.Range("D33").Value = "'0.0913"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.42%  "

$ws.Range("D34").Value = "'5.49"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.15%  "

$ws.Range("E35").Value = "  +2.18%  "

$ws.Range("D36").Value = "'0.0364"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.10%  "

$ws.Range("E37").Value = "  -3.36%  "

$ws.Range("D38").Value = "'4.45"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.63%  "

$ws.Range("D39").Value = "'3.35"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.31%  "

$ws.Range("E40").Value = "  +8.09%  "

$ws.Range("E41").Value = "  +10.74%  "

$ws.Range("D42").Value = "'1.42"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +22.09%  "

$ws.Range("D43").Value = "'12.36"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.12%  "

$ws.Range("D44").Value = "'62.90"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.05%  "

$ws.Range("D45").Value = "'9.14"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +9.00%  "

$ws.Range("E46").Value = "  -1.26%  "

$ws.Range("E47").Value = "  +4.08%  "

$ws.Range("D48").Value = "'100.49"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.58%  "

$ws.Range("D49").Value = "'1.20"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.11%  "

$ws.Range("E50").Value = "  +15.53%  "

$ws.Range("D51").Value = "'2.542.29"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.21%  "
